$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value = 2897.3  # H19: 3250 -> 2897.3
$ws.Cells.Item(19, 9).Value = 3374  # I19: 4249.5 -> 3374
$ws.Cells.Item(19, 10).Value = 2579.5  # J19: 1917.3334 -> 2579.5
$ws.Cells.Item(19, 11).Value = 3374  # K19: 4249.5 -> 3374
$ws.Cells.Item(19, 12).Value = 2579.5  # L19: 1917.3334 -> 2579.5
$ws.Cells.Item(19, 13).Value = -3199  # M19: -4074.5 -> -3199
$ws.Cells.Item(19, 14).Value = -2929.5  # N19: -2267.3334 -> -2929.5

# Row 94
$ws.Cells.Item(94, 8).Value = 3002.5  # H94: 3000 -> 3002.5
$ws.Cells.Item(94, 9).Value = 3002.5  # I94: 3000 -> 3002.5
$ws.Cells.Item(94, 11).Value = 3002.5  # K94: 3000 -> 3002.5
$ws.Cells.Item(94, 13).Value = -2551.5  # M94: -2549 -> -2551.5

# Row 116
$ws.Cells.Item(116, 8).Value = 6533.0835  # H116: 6338.615 -> 6533.0835
$ws.Cells.Item(116, 9).Value = 5858.4287  # I116: 5626.75 -> 5858.4287
$ws.Cells.Item(116, 11).Value = 5858.4287  # K116: 5626.75 -> 5858.4287
$ws.Cells.Item(116, 13).Value = -2416.4287  # M116: -2184.75 -> -2416.4287

# Row 132
$ws.Cells.Item(132, 8).Value = 1204.25  # H132: 1223.2941 -> 1204.25
$ws.Cells.Item(132, 9).Value = 1204.25  # I132: 1223.2941 -> 1204.25
$ws.Cells.Item(132, 11).Value = 3612.75  # K132: 3669.8823 -> 3612.75
$ws.Cells.Item(132, 13).Value = -1082.75  # M132: -1139.8823 -> -1082.75

# Row 135
$ws.Cells.Item(135, 8).Value = 9032.789000000001  # H135: 7985.5 -> 9032.789000000001
$ws.Cells.Item(135, 9).Value = 3601.8572  # I135: 3237.625 -> 3601.8572
$ws.Cells.Item(135, 10).Value = 12200.833  # J135: 10698.571 -> 12200.833
$ws.Cells.Item(135, 11).Value = 32416.7148  # K135: 29138.625 -> 32416.7148
$ws.Cells.Item(135, 12).Value = 109807.497  # L135: 96287.139 -> 109807.497
$ws.Cells.Item(135, 13).Value = -29881.7148  # M135: -26603.625 -> -29881.7148
$ws.Cells.Item(135, 14).Value = -114877.497  # N135: -101357.139 -> -114877.497

# Row 138
$ws.Cells.Item(138, 8).Value = 905.2727  # H138: 1992.375 -> 905.2727
$ws.Cells.Item(138, 9).Value = 905.2727  # I138: 937.8 -> 905.2727
$ws.Cells.Item(138, 10).Value = 0  # J138: 3750 -> 0
$ws.Cells.Item(138, 11).Value = 2715.8181  # K138: 2813.4 -> 2715.8181
$ws.Cells.Item(138, 12).Value = 0  # L138: 11250 -> 0
$ws.Cells.Item(138, 13).Value = 2424.1819  # M138: 2326.6 -> 2424.1819
$ws.Cells.Item(138, 14).Value = $null  # N138: clear (was -21530)

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 19233158  # H32: 19234666 -> 19233158
$ws.Cells.Item(32, 9).Value = 20835172  # I32: 20836804 -> 20835172
$ws.Cells.Item(32, 11).Value = 20835172  # K32: 20836804 -> 20835172
$ws.Cells.Item(32, 13).Value = -20834885  # M32: -20836517 -> -20834885

# Row 34
$ws.Cells.Item(34, 8).Value = 500000  # H34: 265247.5 -> 500000
$ws.Cells.Item(34, 10).Value = 500000  # J34: 265247.5 -> 500000
$ws.Cells.Item(34, 12).Value = 500000  # L34: 265247.5 -> 500000
$ws.Cells.Item(34, 14).Value = -500542  # N34: -265789.5 -> -500542

# Row 74
$ws.Cells.Item(74, 8).Value = 8629785  # H74: 10010288 -> 8629785
$ws.Cells.Item(74, 9).Value = 13891072  # I74: 17859482 -> 13891072
$ws.Cells.Item(74, 11).Value = 13891072  # K74: 17859482 -> 13891072
$ws.Cells.Item(74, 13).Value = -13890198  # M74: -17858608 -> -13890198

# Row 77
$ws.Cells.Item(77, 8).Value = 8629785  # H77: 10010288 -> 8629785
$ws.Cells.Item(77, 9).Value = 13891072  # I77: 17859482 -> 13891072
$ws.Cells.Item(77, 11).Value = 69455360  # K77: 89297410 -> 69455360
$ws.Cells.Item(77, 13).Value = -69450992  # M77: -89293042 -> -69450992

# Row 95
$ws.Cells.Item(95, 8).Value = 39995  # H95: 39990 -> 39995
$ws.Cells.Item(95, 10).Value = 39995  # J95: 39990 -> 39995
$ws.Cells.Item(95, 12).Value = 39995  # L95: 39990 -> 39995
$ws.Cells.Item(95, 14).Value = -45487  # N95: -45482 -> -45487

# Row 122
$ws.Cells.Item(122, 8).Value = 2425.625  # H122: 2533.9333 -> 2425.625
$ws.Cells.Item(122, 9).Value = 1351.5  # I122: 1502 -> 1351.5
$ws.Cells.Item(122, 10).Value = 3499.75  # J122: 3221.889 -> 3499.75
$ws.Cells.Item(122, 11).Value = 4054.5  # K122: 4506 -> 4054.5
$ws.Cells.Item(122, 12).Value = 10499.25  # L122: 9665.667000000001 -> 10499.25
$ws.Cells.Item(122, 13).Value = -1604.5  # M122: -2056 -> -1604.5
$ws.Cells.Item(122, 14).Value = -15399.25  # N122: -14565.667 -> -15399.25

# Row 132
$ws.Cells.Item(132, 8).Value = 4670.205  # H132: 5106.514 -> 4670.205
$ws.Cells.Item(132, 9).Value = 2488.0386  # I132: 2785.4092 -> 2488.0386
$ws.Cells.Item(132, 11).Value = 7464.1158  # K132: 8356.2276 -> 7464.1158
$ws.Cells.Item(132, 13).Value = -4934.1158  # M132: -5826.2276 -> -4934.1158

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 3910.55  # H20: 4037 -> 3910.55
$ws.Cells.Item(20, 9).Value = 3938.5  # I20: 4081.4707 -> 3938.5
$ws.Cells.Item(20, 11).Value = 3938.5  # K20: 4081.4707 -> 3938.5
$ws.Cells.Item(20, 13).Value = -3691.5  # M20: -3834.4707 -> -3691.5

# Row 86
$ws.Cells.Item(86, 8).Value = 1723.762  # H86: 1766.0476 -> 1723.762
$ws.Cells.Item(86, 9).Value = 1695  # I86: 1753.1578 -> 1695
$ws.Cells.Item(86, 10).Value = 1896.3334  # J86: 1888.5 -> 1896.3334
$ws.Cells.Item(86, 11).Value = 1695  # K86: 1753.1578 -> 1695
$ws.Cells.Item(86, 12).Value = 1896.3334  # L86: 1888.5 -> 1896.3334
$ws.Cells.Item(86, 13).Value = -572  # M86: -630.1578 -> -572
$ws.Cells.Item(86, 14).Value = -4142.3334  # N86: -4134.5 -> -4142.3334

# Row 89
$ws.Cells.Item(89, 8).Value = 1723.762  # H89: 1766.0476 -> 1723.762
$ws.Cells.Item(89, 9).Value = 1695  # I89: 1753.1578 -> 1695
$ws.Cells.Item(89, 10).Value = 1896.3334  # J89: 1888.5 -> 1896.3334
$ws.Cells.Item(89, 11).Value = 8475  # K89: 8765.789000000001 -> 8475
$ws.Cells.Item(89, 12).Value = 9481.666999999999  # L89: 9442.5 -> 9481.666999999999
$ws.Cells.Item(89, 13).Value = -2859  # M89: -3149.789000000001 -> -2859
$ws.Cells.Item(89, 14).Value = -20713.667  # N89: -20674.5 -> -20713.667

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 1494.0741  # H7: 1495.8148 -> 1494.0741
$ws.Cells.Item(7, 9).Value = 134.4  # I7: 132.0625 -> 134.4
$ws.Cells.Item(7, 10).Value = 3193.6667  # J7: 3479.4546 -> 3193.6667
$ws.Cells.Item(7, 11).Value = 134.4  # K7: 132.0625 -> 134.4
$ws.Cells.Item(7, 12).Value = 3193.6667  # L7: 3479.4546 -> 3193.6667
$ws.Cells.Item(7, 13).Value = -21.40000000000001  # M7: -19.0625 -> -21.40000000000001
$ws.Cells.Item(7, 14).Value = -3419.6667  # N7: -3705.4546 -> -3419.6667

# Row 22
$ws.Cells.Item(22, 8).Value = 674  # H22: 475.5 -> 674
$ws.Cells.Item(22, 9).Value = 401  # I22: 389.75 -> 401
$ws.Cells.Item(22, 10).Value = 947  # J22: 561.25 -> 947
$ws.Cells.Item(22, 11).Value = 401  # K22: 389.75 -> 401
$ws.Cells.Item(22, 12).Value = 947  # L22: 561.25 -> 947
$ws.Cells.Item(22, 13).Value = -51  # M22: -39.75 -> -51
$ws.Cells.Item(22, 14).Value = -1647  # N22: -1261.25 -> -1647

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 3187.125  # H131: 3332.9333 -> 3187.125
$ws.Cells.Item(131, 9).Value = 1292  # I131: 1299.6364 -> 1292
$ws.Cells.Item(131, 10).Value = 11399.333  # J131: 8924.5 -> 11399.333
$ws.Cells.Item(131, 11).Value = 3876  # K131: 3898.9092 -> 3876
$ws.Cells.Item(131, 12).Value = 34197.999  # L131: 26773.5 -> 34197.999
$ws.Cells.Item(131, 13).Value = 1164  # M131: 1141.0908 -> 1164
$ws.Cells.Item(131, 14).Value = -44277.999  # N131: -36853.5 -> -44277.999

# Row 133
$ws.Cells.Item(133, 8).Value = 3424.375  # H133: 3428.125 -> 3424.375
$ws.Cells.Item(133, 9).Value = 3424.375  # I133: 3428.125 -> 3424.375
$ws.Cells.Item(133, 11).Value = 10273.125  # K133: 10284.375 -> 10273.125
$ws.Cells.Item(133, 13).Value = -5213.125  # M133: -5224.375 -> -5213.125

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Cells.Item(11, 8).Value = 30034286  # H11: 17199858 -> 30034286
$ws.Cells.Item(11, 9).Value = 40030000  # I11: 30037500 -> 40030000
$ws.Cells.Item(11, 10).Value = 5045002  # J11: 83002.664 -> 5045002
$ws.Cells.Item(11, 11).Value = 40030000  # K11: 30037500 -> 40030000
$ws.Cells.Item(11, 12).Value = 5045002  # L11: 83002.664 -> 5045002
$ws.Cells.Item(11, 13).Value = -40029861  # M11: -30037361 -> -40029861
$ws.Cells.Item(11, 14).Value = -5045280  # N11: -83280.664 -> -5045280

# Row 18
$ws.Cells.Item(18, 8).Value = 30000  # H18: 29995 -> 30000
$ws.Cells.Item(18, 9).Value = 30000  # I18: 29995 -> 30000
$ws.Cells.Item(18, 11).Value = 30000  # K18: 29995 -> 30000
$ws.Cells.Item(18, 13).Value = -29707  # M18: -29702 -> -29707

# Row 20
$ws.Cells.Item(20, 8).Value = 30548.924  # H20: 27769.732 -> 30548.924
$ws.Cells.Item(20, 9).Value = 10000  # I20: 9800 -> 10000
$ws.Cells.Item(20, 10).Value = 32261.334  # J20: 30534.309 -> 32261.334
$ws.Cells.Item(20, 11).Value = 10000  # K20: 9800 -> 10000
$ws.Cells.Item(20, 12).Value = 32261.334  # L20: 30534.309 -> 32261.334
$ws.Cells.Item(20, 13).Value = -9755  # M20: -9555 -> -9755
$ws.Cells.Item(20, 14).Value = -32751.334  # N20: -31024.309 -> -32751.334

# Row 46
$ws.Cells.Item(46, 8).Value = 0  # H46: 5870 -> 0
$ws.Cells.Item(46, 10).Value = 0  # J46: 5870 -> 0
$ws.Cells.Item(46, 12).Value = 0  # L46: 5870 -> 0
$ws.Cells.Item(46, 14).Value = $null  # N46: clear (was -6182)

# Row 70
$ws.Cells.Item(70, 8).Value = 7480.0586  # H70: 7680.1177 -> 7480.0586
$ws.Cells.Item(70, 9).Value = 6285.2856  # I70: 6771.143 -> 6285.2856
$ws.Cells.Item(70, 11).Value = 6285.2856  # K70: 6771.143 -> 6285.2856
$ws.Cells.Item(70, 13).Value = -6015.2856  # M70: -6501.143 -> -6015.2856

# Row 73
$ws.Cells.Item(73, 8).Value = 7480.0586  # H73: 7680.1177 -> 7480.0586
$ws.Cells.Item(73, 9).Value = 6285.2856  # I73: 6771.143 -> 6285.2856
$ws.Cells.Item(73, 11).Value = 6285.2856  # K73: 6771.143 -> 6285.2856
$ws.Cells.Item(73, 13).Value = -5349.2856  # M73: -5835.143 -> -5349.2856

# Row 97
$ws.Cells.Item(97, 8).Value = 880.0476  # H97: 922.1579 -> 880.0476
$ws.Cells.Item(97, 9).Value = 693.4  # I97: 690 -> 693.4
$ws.Cells.Item(97, 10).Value = 1346.6666  # J97: 1792.75 -> 1346.6666
$ws.Cells.Item(97, 11).Value = 693.4  # K97: 690 -> 693.4
$ws.Cells.Item(97, 12).Value = 1346.6666  # L97: 1792.75 -> 1346.6666
$ws.Cells.Item(97, 13).Value = -197.4  # M97: -194 -> -197.4
$ws.Cells.Item(97, 14).Value = -2338.6666  # N97: -2784.75 -> -2338.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Cells.Item(20, 8).Value = 30000  # H20: 9780.6 -> 30000
$ws.Cells.Item(20, 9).Value = 0  # I20: 5400 -> 0
$ws.Cells.Item(20, 10).Value = 30000  # J20: 10875.75 -> 30000
$ws.Cells.Item(20, 11).Value = 0  # K20: 5400 -> 0
$ws.Cells.Item(20, 12).Value = 30000  # L20: 10875.75 -> 30000
$ws.Cells.Item(20, 13).Value = $null  # M20: clear (was -5174)
$ws.Cells.Item(20, 14).Value = -30452  # N20: -11327.75 -> -30452

# Row 46
$ws.Cells.Item(46, 8).Value = 3349.7083  # H46: 3295.6365 -> 3349.7083
$ws.Cells.Item(46, 9).Value = 3268.6924  # I46: 3145.818 -> 3268.6924
$ws.Cells.Item(46, 11).Value = 3268.6924  # K46: 3145.818 -> 3268.6924
$ws.Cells.Item(46, 13).Value = -3080.6924  # M46: -2957.818 -> -3080.6924

# Row 82
$ws.Cells.Item(82, 8).Value = 2000  # H82: 2220.2 -> 2000
$ws.Cells.Item(82, 10).Value = 1949.5  # J82: 2299.6667 -> 1949.5
$ws.Cells.Item(82, 12).Value = 1949.5  # L82: 2299.6667 -> 1949.5
$ws.Cells.Item(82, 14).Value = -2671.5  # N82: -3021.6667 -> -2671.5

# Row 85
$ws.Cells.Item(85, 8).Value = 2000  # H85: 2220.2 -> 2000
$ws.Cells.Item(85, 10).Value = 1949.5  # J85: 2299.6667 -> 1949.5
$ws.Cells.Item(85, 12).Value = 1949.5  # L85: 2299.6667 -> 1949.5
$ws.Cells.Item(85, 14).Value = -4445.5  # N85: -4795.6667 -> -4445.5

# Row 93
$ws.Cells.Item(93, 8).Value = 33334334  # H93: 40001056 -> 33334334
$ws.Cells.Item(93, 9).Value = 40000970  # I93: 45455564 -> 40000970
$ws.Cells.Item(93, 10).Value = 1169.2  # J93: 1332.6666 -> 1169.2
$ws.Cells.Item(93, 11).Value = 40000970  # K93: 45455564 -> 40000970
$ws.Cells.Item(93, 12).Value = 1169.2  # L93: 1332.6666 -> 1169.2
$ws.Cells.Item(93, 13).Value = -39999722  # M93: -45454316 -> -39999722
$ws.Cells.Item(93, 14).Value = -3665.2  # N93: -3828.6666 -> -3665.2

# Row 138
$ws.Cells.Item(138, 8).Value = 60000  # H138: 0 -> 60000
$ws.Cells.Item(138, 10).Value = 60000  # J138: 0 -> 60000
$ws.Cells.Item(138, 12).Value = 60000  # L138: 0 -> 60000
$ws.Cells.Item(138, 14).Value = -70280  # N138: None -> -70280

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Cells.Item(75, 8).Value = 12543308  # H75: 15044315 -> 12543308
$ws.Cells.Item(75, 10).Value = 12543308  # J75: 15044315 -> 12543308
$ws.Cells.Item(75, 12).Value = 12543308  # L75: 15044315 -> 12543308
$ws.Cells.Item(75, 14).Value = -12545180  # N75: -15046187 -> -12545180

# Row 78
$ws.Cells.Item(78, 8).Value = 12543308  # H78: 15044315 -> 12543308
$ws.Cells.Item(78, 10).Value = 12543308  # J78: 15044315 -> 12543308
$ws.Cells.Item(78, 12).Value = 37629924  # L78: 45132945 -> 37629924
$ws.Cells.Item(78, 14).Value = -37639284  # N78: -45142305 -> -37639284

# Row 116
$ws.Cells.Item(116, 8).Value = 96990  # H116: 0 -> 96990
$ws.Cells.Item(116, 10).Value = 96990  # J116: 0 -> 96990
$ws.Cells.Item(116, 12).Value = 96990  # L116: 0 -> 96990
$ws.Cells.Item(116, 14).Value = -106168  # N116: None -> -106168

# Row 122
$ws.Cells.Item(122, 8).Value = 4384.875  # H122: 4597.567 -> 4384.875
$ws.Cells.Item(122, 10).Value = 6953.273  # J122: 8233 -> 6953.273
$ws.Cells.Item(122, 12).Value = 20859.819  # L122: 24699 -> 20859.819
$ws.Cells.Item(122, 14).Value = -25759.819  # N122: -29599 -> -25759.819

# Row 133
$ws.Cells.Item(133, 8).Value = 98922.5  # H133: 99125 -> 98922.5
$ws.Cells.Item(133, 10).Value = 98922.5  # J133: 99125 -> 98922.5
$ws.Cells.Item(133, 12).Value = 98922.5  # L133: 99125 -> 98922.5
$ws.Cells.Item(133, 14).Value = -109042.5  # N133: -109245 -> -109042.5

# Row 136
$ws.Cells.Item(136, 8).Value = 1333.0526  # H136: 1382.4166 -> 1333.0526
$ws.Cells.Item(136, 9).Value = 1024.5172  # I136: 1067.4814 -> 1024.5172
$ws.Cells.Item(136, 11).Value = 3073.5516  # K136: 3202.4442 -> 3073.5516
$ws.Cells.Item(136, 13).Value = -523.5515999999998  # M136: -652.4441999999999 -> -523.5515999999998
